$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prompt = "`nYou are a change management expert helping companies adopt new technologies using structured frameworks.`nA company wants to adopt a new technology: grok`nThey want to follow the adkar framework.`nGive a clear, step-by-step guide tailored to ops. using the adkar model. Make it practical and actionable.`nUse numbered steps and keep the tone helpful.`nAlso include a list of common technical FAQs regarding the technology and answers to help employees understand the change process.`n"

# Row 4
$ws.Range("A4").Value = "689580fc-3200-4fe8-a3b6-02eb7cb60de4"
$ws.Range("B4").Value = "2025-03-23 00:41:49"
$ws.Range("C4").Value = $prompt
$ws.Range("D4").Value = "make it easier for ops"

# Row 5
$ws.Range("A5").Value = "4a1f6efc-bf0c-40b3-b27e-cfcb56cc26e5"
$ws.Range("B5").Value = "2025-03-23 00:43:42"
$ws.Range("C5").Value = $prompt
$ws.Range("D5").Value = "use levin"
